$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 124: new time-record entry
$ws.Cells.Item(124, 1).Value = 44619
$ws.Cells.Item(124, 2).Value = "Editing puzzles, enviroment and erd"
$ws.Cells.Item(124, 3).NumberFormat = "0.00"
$ws.Cells.Item(124, 3).Value = 1.42
$ws.Cells.Item(124, 4).Value = "Sarvan Amel"

# Update total for "Lasinger Christoph" (row 129) - B130 formula recalculates automatically
$ws.Cells.Item(129, 2).Value = 15

# Update the current selection to match the saved workbook state
$ws.Range("D125").Select()
